$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (sheet ALC)
$ws.Range("H33").Value = 221.625
$ws.Range("I33").Value = 221.625
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 221.625
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 7.375
$ws.Range("N33").ClearContents()

# Row 58 (sheet ALC)
$ws.Range("H58").Value = 7477.125
$ws.Range("I58").Value = 450
$ws.Range("J58").Value = 9819.5
$ws.Range("K58").Value = 1350
$ws.Range("L58").Value = 29458.5
$ws.Range("M58").Value = -1200
$ws.Range("N58").Value = -29758.5

# Row 111 (sheet ALC)
$ws.Range("H111").Value = 150
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 150
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 450
$ws.Range("N111").Value = -6584

# Row 113 (sheet ALC)
$ws.Range("H113").Value = 1114222.2
$ws.Range("I113").Value = 10000000

# Row 115 (sheet ALC)
$ws.Range("H115").Value = 15000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 15000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 45000
$ws.Range("N115").Value = -48134
$ws.Range("M115").ClearContents()

# Row 116 (sheet ALC)
$ws.Range("H116").Value = 9997.25
$ws.Range("I116").Value = 9996.333000000001
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 9996.333000000001
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -6554.333000000001
$ws.Range("N116").Value = -16884

# Row 118 (sheet ALC)
$ws.Range("H118").Value = 2799.6667
$ws.Range("I118").Value = 2800
$ws.Range("J118").Value = 2799.5
$ws.Range("K118").Value = 8400
$ws.Range("L118").Value = 8398.5
$ws.Range("M118").Value = -6743
$ws.Range("N118").Value = -11712.5

# Row 125 (sheet ALC)
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -15540

# Row 135 (sheet ALC)
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 102 (sheet ARM)
$ws.Range("H102").Value = 7450
$ws.Range("I102").Value = 7450
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 7450
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -5828

# Row 122 (sheet ARM)
$ws.Range("H122").Value = 5753.1665
$ws.Range("I122").Value = 6103.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 18311.4
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -15861.4
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (sheet BSM)
$ws.Range("H86").Value = 2488.5454
$ws.Range("I86").Value = 2488.5454
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2488.5454
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1365.5454

# Row 89 (sheet BSM)
$ws.Range("H89").Value = 2488.5454
$ws.Range("I89").Value = 2488.5454
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12442.727
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6826.726999999999

# Row 94 (sheet BSM)
$ws.Range("H94").Value = 1367.3334
$ws.Range("I94").Value = 974.2
$ws.Range("J94").Value = 3333
$ws.Range("K94").Value = 974.2
$ws.Range("L94").Value = 3333
$ws.Range("M94").Value = -523.2
$ws.Range("N94").Value = -4235

# Row 105 (sheet BSM)
$ws.Range("H105").Value = 5412.6665
$ws.Range("I105").Value = 4995.2
$ws.Range("J105").Value = 7500
$ws.Range("K105").Value = 4995.2
$ws.Range("L105").Value = 7500
$ws.Range("M105").Value = -3248.2
$ws.Range("N105").Value = -10994

# Row 134 (sheet BSM)
$ws.Range("H134").Value = 2689.7144
$ws.Range("I134").Value = 2650.5454
$ws.Range("J134").Value = 2833.3333
$ws.Range("K134").Value = 7951.6362
$ws.Range("L134").Value = 8499.999899999999
$ws.Range("M134").Value = -5416.6362
$ws.Range("N134").Value = -13569.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (sheet CRP)
$ws.Range("H31").Value = 2339.8462
$ws.Range("I31").Value = 2615.7144
$ws.Range("J31").Value = 2018
$ws.Range("K31").Value = 2615.7144
$ws.Range("L31").Value = 2018
$ws.Range("M31").Value = -2320.7144
$ws.Range("N31").Value = -2608

# Row 34 (sheet CRP)
$ws.Range("H34").Value = 2339.8462
$ws.Range("I34").Value = 2615.7144
$ws.Range("J34").Value = 2018
$ws.Range("K34").Value = 2615.7144
$ws.Range("L34").Value = 2018
$ws.Range("M34").Value = -2413.7144
$ws.Range("N34").Value = -2422

# Row 99 (sheet CRP)
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 10000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -8502

# Row 122 (sheet CRP)
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("M122").ClearContents()

# Row 126 (sheet CRP)
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 30000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27530

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (sheet CUL)
$ws.Range("H34").Value = 4349.75
$ws.Range("I34").Value = 1366
$ws.Range("J34").Value = 6140
$ws.Range("K34").Value = 4098
$ws.Range("L34").Value = 18420
$ws.Range("M34").Value = -4014
$ws.Range("N34").Value = -18588

# Row 87 (sheet CUL)
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()

# Row 90 (sheet CUL)
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()

# Row 113 (sheet CUL)
$ws.Range("H113").Value = 2682.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2682.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8047.5
$ws.Range("N113").Value = -12387.5

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (sheet GSM)
$ws.Range("H113").Value = 1816.6666
$ws.Range("I113").Value = 1816.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1816.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 353.3334

# Row 126 (sheet GSM)
$ws.Range("H126").Value = 9999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 29997
$ws.Range("N126").Value = -34937

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (sheet LTW)
$ws.Range("H68").Value = 2666
$ws.Range("I68").Value = 2999
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2999
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -2250
$ws.Range("N68").Value = -3498

# Row 71 (sheet LTW)
$ws.Range("H71").Value = 2666
$ws.Range("I71").Value = 2999
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 14995
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -11251
$ws.Range("N71").Value = -17488

# Row 122 (sheet LTW)
$ws.Range("H122").Value = 10160
$ws.Range("I122").Value = 6800
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 20400
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -17950
$ws.Range("N122").Value = -37900

# Row 132 (sheet LTW)
$ws.Range("H132").Value = 5474.5
$ws.Range("I132").Value = 4950
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 14850
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -12320
$ws.Range("N132").Value = -23057

# Row 136 (sheet LTW)
$ws.Range("H136").Value = 5196.8887
$ws.Range("I136").Value = 6414.6665
$ws.Range("J136").Value = 2761.3333
$ws.Range("K136").Value = 19243.9995
$ws.Range("L136").Value = 8283.999899999999
$ws.Range("M136").Value = -16693.9995
$ws.Range("N136").Value = -13383.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (sheet WVR)
$ws.Range("H107").Value = 535.2727
$ws.Range("I107").Value = 509.77777
$ws.Range("J107").Value = 650
$ws.Range("K107").Value = 1529.33331
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = 390.66669
$ws.Range("N107").Value = -5790

# Row 122 (sheet WVR)
$ws.Range("H122").Value = 3965.889
$ws.Range("I122").Value = 3965.889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11897.667
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9447.667000000001

# Row 126 (sheet WVR)
$ws.Range("H126").Value = 1582.0769
$ws.Range("I126").Value = 1582.0769
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4746.2307
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2276.2307

# Row 136 (sheet WVR)
$ws.Range("H136").Value = 3207.2778
$ws.Range("I136").Value = 3162.5
$ws.Range("J136").Value = 3296.8333
$ws.Range("K136").Value = 9487.5
$ws.Range("L136").Value = 9890.499899999999
$ws.Range("M136").Value = -6937.5
$ws.Range("N136").Value = -14990.4999
